$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change C2 from "tes@tes" to "tes@te"
$ws.Range("C2").Value = "tes@te"

# Fill in C3 with new value "apapun@telkom.co.id"
$ws.Range("C3").Value = "apapun@telkom.co.id"
